$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 456, shifting existing rows 456:553 down to 457:554.
$ws.Range("A456:R456").EntireRow.Insert()

# Populate the newly inserted row 456 with a copy of the (now shifted) row 457's
# surrounding context plus the new reading's own measurements.
$ws.Range("A456").Value = 10
$ws.Range("B456").Value = "Vega Modelo de Temuco"
$ws.Range("C456").Value = "La Araucanía"
$ws.Range("D456").Value = 44798
$ws.Range("E456").Value = 9
$ws.Range("F456").Value = 100112032
$ws.Range("G456").Value = "Zapallo italiano"
$ws.Range("H456").Value = "Sin especificar"
$ws.Range("I456").Value = "Primera"
$ws.Range("J456").Value = 500
$ws.Range("K456").Value = 24000
$ws.Range("L456").Value = 25000
$ws.Range("M456").Value = 24400
$ws.Range("N456").Value = "$/caja 60 unidades"
$ws.Range("O456").Value = "Región de Arica y Parinacota"
$ws.Range("P456").Value = 407
$ws.Range("Q456").Value = 60
$ws.Range("R456").Value = "Hortaliza"
